# Consolidate the "<word> " + "<word>" run pairs into a single run,
# leaving the trailing word run untouched - mirrors the writer change
# that merges adjacent text runs to shrink the generated file.

$p = $ppt.ActivePresentation

# Slide 2: TextBox "The" + " " + "Moon"  ->  "The " + "Moon"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Characters(1, 4).Text = "The "

# Slide 3: Title "One" + " " + "More"  ->  "One " + "More"
$s3 = $p.Slides.Item(3)
$tr3title = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3title.Characters(1, 4).Text = "One "

# Slide 3: TextBox "The" + " " + "Moon"  ->  "The " + "Moon"
$tr3box = $s3.Shapes.Item(3).TextFrame.TextRange
$tr3box.Characters(1, 4).Text = "The "
